# This edit inserts 3 new weekly price records (rows 254-256) for
# "Tuna" at "Provincia de Chacabuco" dated 2022-03-24, pushing the
# previously-existing rows 254-316 down to 257-319.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 254; this shifts rows 254:316 down to
# 257:319 automatically (Excel also extends the used range / dimension).
$ws.Rows("254:256").Insert()

# Common (constant-across-the-sheet) column values for this product.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107011
$categoria = "Tuna"
$variedad  = "Sin especificar"

# New row 254: Especial
$ws.Cells.Item(254, 1).Value  = $mercadoId
$ws.Cells.Item(254, 2).Value  = $mercado
$ws.Cells.Item(254, 3).Value  = $region
$ws.Cells.Item(254, 4).Value  = 44644
$ws.Cells.Item(254, 5).Value  = $codreg
$ws.Cells.Item(254, 6).Value  = $tipo
$ws.Cells.Item(254, 7).Value  = $productoId
$ws.Cells.Item(254, 8).Value  = $producto
$ws.Cells.Item(254, 9).Value  = $categoriaId
$ws.Cells.Item(254, 10).Value = $categoria
$ws.Cells.Item(254, 11).Value = $variedad
$ws.Cells.Item(254, 12).Value = "Especial"
$ws.Cells.Item(254, 13).Value = 180
$ws.Cells.Item(254, 14).Value = 15000
$ws.Cells.Item(254, 15).Value = 15000
$ws.Cells.Item(254, 16).Value = 15000
$ws.Cells.Item(254, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(254, 18).Value = "Provincia de Chacabuco"
$ws.Cells.Item(254, 19).Value = 833
$ws.Cells.Item(254, 20).Value = 18

# New row 255: Primera
$ws.Cells.Item(255, 1).Value  = $mercadoId
$ws.Cells.Item(255, 2).Value  = $mercado
$ws.Cells.Item(255, 3).Value  = $region
$ws.Cells.Item(255, 4).Value  = 44644
$ws.Cells.Item(255, 5).Value  = $codreg
$ws.Cells.Item(255, 6).Value  = $tipo
$ws.Cells.Item(255, 7).Value  = $productoId
$ws.Cells.Item(255, 8).Value  = $producto
$ws.Cells.Item(255, 9).Value  = $categoriaId
$ws.Cells.Item(255, 10).Value = $categoria
$ws.Cells.Item(255, 11).Value = $variedad
$ws.Cells.Item(255, 12).Value = "Primera"
$ws.Cells.Item(255, 13).Value = 200
$ws.Cells.Item(255, 14).Value = 13000
$ws.Cells.Item(255, 15).Value = 13000
$ws.Cells.Item(255, 16).Value = 13000
$ws.Cells.Item(255, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(255, 18).Value = "Provincia de Chacabuco"
$ws.Cells.Item(255, 19).Value = 722
$ws.Cells.Item(255, 20).Value = 18

# New row 256: Segunda
$ws.Cells.Item(256, 1).Value  = $mercadoId
$ws.Cells.Item(256, 2).Value  = $mercado
$ws.Cells.Item(256, 3).Value  = $region
$ws.Cells.Item(256, 4).Value  = 44644
$ws.Cells.Item(256, 5).Value  = $codreg
$ws.Cells.Item(256, 6).Value  = $tipo
$ws.Cells.Item(256, 7).Value  = $productoId
$ws.Cells.Item(256, 8).Value  = $producto
$ws.Cells.Item(256, 9).Value  = $categoriaId
$ws.Cells.Item(256, 10).Value = $categoria
$ws.Cells.Item(256, 11).Value = $variedad
$ws.Cells.Item(256, 12).Value = "Segunda"
$ws.Cells.Item(256, 13).Value = 150
$ws.Cells.Item(256, 14).Value = 8000
$ws.Cells.Item(256, 15).Value = 8000
$ws.Cells.Item(256, 16).Value = 8000
$ws.Cells.Item(256, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(256, 18).Value = "Provincia de Chacabuco"
$ws.Cells.Item(256, 19).Value = 444
$ws.Cells.Item(256, 20).Value = 18

# Ensure the date cells keep the same date number format used by the
# rest of the "Fecha" column.
$ws.Range("D254:D256").NumberFormat = $ws.Range("D257").NumberFormat
